$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.642.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.163.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.90%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.04%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.165.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.67%  "

$ws.Range("E10").Value = "  -4.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("E12").Value = "  -4.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.680.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.761.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.83%  "

$ws.Range("E17").Value = "  +0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.166.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.11%  "

$ws.Range("E19").Value = "  -3.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.33%  "

$ws.Range("E21").Value = "  -1.60%  "

$ws.Range("E22").Value = "  -1.87%  "

$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  -2.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.73%  "

$ws.Range("E29").Value = "  -3.72%  "

$ws.Range("E30").Value = "  +2.40%  "

$ws.Range("E31").Value = "  +2.29%  "

$ws.Range("E32").Value = "  -6.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.20%  "

$ws.Range("E35").Value = "  -5.33%  "

$ws.Range("E36").Value = "  -4.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0750"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.66%  "

$ws.Range("E41").Value = "  -3.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0403"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.24%  "

$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.901.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("E46").Value = "  -6.87%  "

$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("E49").Value = "  -2.27%  "

$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("E51").Value = "  -1.95%  "
